$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1654.7142
$ws.Range("I43").Value = 1556.2
$ws.Range("J43").Value = 1901
$ws.Range("K43").Value = 1556.2
$ws.Range("L43").Value = 1901
$ws.Range("M43").Value = -1487.2
$ws.Range("N43").Value = -2039

$ws.Range("H98").Value = 49104.332
$ws.Range("I98").Value = 63724.438
$ws.Range("K98").Value = 63724.438
$ws.Range("M98").Value = -62226.438

$ws.Range("H122").Value = 49104.332
$ws.Range("I122").Value = 63724.438
$ws.Range("K122").Value = 191173.314
$ws.Range("M122").Value = -188723.314

$ws.Range("H132").Value = 1431.36
$ws.Range("I132").Value = 990.6087
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 2971.8261
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -441.8261000000002
$ws.Range("N132").Value = -24560

$ws.Range("H137").Value = 1055
$ws.Range("I137").Value = 972.6786
$ws.Range("J137").Value = 1516
$ws.Range("K137").Value = 2918.0358
$ws.Range("L137").Value = 4548
$ws.Range("M137").Value = -368.0357999999997
$ws.Range("N137").Value = -9648

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1338.5834
$ws.Range("I110").Value = 922.4211
$ws.Range("J110").Value = 2920
$ws.Range("K110").Value = 922.4211
$ws.Range("L110").Value = 2920
$ws.Range("M110").Value = 1122.5789
$ws.Range("N110").Value = -7010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H99").Value = 1936
$ws.Range("I99").Value = 2260
$ws.Range("J99").Value = 1450
$ws.Range("K99").Value = 2260
$ws.Range("L99").Value = 1450
$ws.Range("M99").Value = -762
$ws.Range("N99").Value = -4446

$ws.Range("H134").Value = 2615.5833
$ws.Range("I134").Value = 2487.375
$ws.Range("J134").Value = 2872
$ws.Range("K134").Value = 7462.125
$ws.Range("L134").Value = 8616
$ws.Range("M134").Value = -4927.125
$ws.Range("N134").Value = -13686

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5586.625
$ws.Range("I16").Value = 3720
$ws.Range("J16").Value = 7453.25
$ws.Range("K16").Value = 3720
$ws.Range("L16").Value = 7453.25
$ws.Range("M16").Value = -3433
$ws.Range("N16").Value = -8027.25

$ws.Range("H22").Value = 100000450
$ws.Range("I22").Value = 500000000
$ws.Range("J22").Value = 562.5
$ws.Range("K22").Value = 500000000
$ws.Range("L22").Value = 562.5
$ws.Range("M22").Value = -499999650
$ws.Range("N22").Value = -1262.5

$ws.Range("H31").Value = 1104.5405
$ws.Range("I31").Value = 851.2646999999999
$ws.Range("J31").Value = 3975
$ws.Range("K31").Value = 851.2646999999999
$ws.Range("L31").Value = 3975
$ws.Range("M31").Value = -556.2646999999999
$ws.Range("N31").Value = -4565

$ws.Range("H34").Value = 1104.5405
$ws.Range("I34").Value = 851.2646999999999
$ws.Range("J34").Value = 3975
$ws.Range("K34").Value = 851.2646999999999
$ws.Range("L34").Value = 3975
$ws.Range("M34").Value = -649.2646999999999
$ws.Range("N34").Value = -4379

$ws.Range("H58").Value = 1705.5
$ws.Range("I58").Value = 411
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 411
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -208
$ws.Range("N58").Value = -3406

$ws.Range("H74").Value = 15322.286
$ws.Range("J74").Value = 15322.286
$ws.Range("L74").Value = 15322.286
$ws.Range("N74").Value = -17070.286

$ws.Range("H77").Value = 15322.286
$ws.Range("J77").Value = 15322.286
$ws.Range("L77").Value = 45966.858
$ws.Range("N77").Value = -54702.858

$ws.Range("H88").Value = 39471.25
$ws.Range("J88").Value = 39471.25
$ws.Range("L88").Value = 39471.25
$ws.Range("N88").Value = -40283.25

$ws.Range("H91").Value = 39471.25
$ws.Range("J91").Value = 39471.25
$ws.Range("L91").Value = 39471.25
$ws.Range("N91").Value = -42279.25

$ws.Range("H107").Value = 855.1429000000001
$ws.Range("I107").Value = 350.47058
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 350.47058
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 1569.52942
$ws.Range("N107").Value = -6840

$ws.Range("H110").Value = 45561.6
$ws.Range("J110").Value = 45561.6
$ws.Range("L110").Value = 45561.6
$ws.Range("N110").Value = -53741.6

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H112").Value = 40234
$ws.Range("J112").Value = 40234
$ws.Range("L112").Value = 40234
$ws.Range("N112").Value = -43188

$ws.Range("H113").Value = 5586.625
$ws.Range("I113").Value = 3720
$ws.Range("J113").Value = 7453.25
$ws.Range("K113").Value = 3720
$ws.Range("L113").Value = 7453.25
$ws.Range("M113").Value = -1550
$ws.Range("N113").Value = -11793.25

$ws.Range("H114").Value = 52228
$ws.Range("J114").Value = 52228
$ws.Range("L114").Value = 52228
$ws.Range("N114").Value = -60906

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 1446.7273
$ws.Range("I122").Value = 1400
$ws.Range("K122").Value = 4200
$ws.Range("M122").Value = -1750

$ws.Range("H132").Value = 2045
$ws.Range("I132").Value = 791.7143
$ws.Range("J132").Value = 3799.6
$ws.Range("K132").Value = 2375.1429
$ws.Range("L132").Value = 11398.8
$ws.Range("M132").Value = 154.8571000000002
$ws.Range("N132").Value = -16458.8

$ws.Range("H134").Value = 1796.7826
$ws.Range("I134").Value = 1341.8823
$ws.Range("J134").Value = 3085.6667
$ws.Range("K134").Value = 4025.6469
$ws.Range("L134").Value = 9257.000100000001
$ws.Range("M134").Value = -1490.6469
$ws.Range("N134").Value = -14327.0001

$ws.Range("H136").Value = 1705.5
$ws.Range("I136").Value = 411
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 1233
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = 1317
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5050.067
$ws.Range("J133").Value = 5671.2856
$ws.Range("L133").Value = 17013.8568
$ws.Range("N133").Value = -27133.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4799.75
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 5199.7144
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 15599.1432
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -20659.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1381.4286
$ws.Range("I16").Value = 1407.5
$ws.Range("J16").Value = 1225
$ws.Range("K16").Value = 1407.5
$ws.Range("L16").Value = 1225
$ws.Range("M16").Value = -1237.5
$ws.Range("N16").Value = -1565

$ws.Range("H55").Value = 140.67392
$ws.Range("I55").Value = 129.79488
$ws.Range("J55").Value = 201.28572
$ws.Range("K55").Value = 129.79488
$ws.Range("L55").Value = 201.28572
$ws.Range("M55").Value = 43.20511999999999
$ws.Range("N55").Value = -547.28572

$ws.Range("H94").Value = 14500
$ws.Range("J94").Value = 14500
$ws.Range("L94").Value = 14500
$ws.Range("N94").Value = -15852

$ws.Range("H132").Value = 2324.122
$ws.Range("I132").Value = 1881.3334
$ws.Range("J132").Value = 3178.0715
$ws.Range("K132").Value = 5644.0002
$ws.Range("L132").Value = 9534.2145
$ws.Range("M132").Value = -3114.0002
$ws.Range("N132").Value = -14594.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 6698.6665
$ws.Range("J41").Value = 6877
$ws.Range("L41").Value = 6877
$ws.Range("N41").Value = -7657

$ws.Range("H132").Value = 1773.1282
$ws.Range("I132").Value = 914.5714
$ws.Range("J132").Value = 2253.92
$ws.Range("K132").Value = 2743.7142
$ws.Range("L132").Value = 6761.76
$ws.Range("M132").Value = -213.7142000000003
$ws.Range("N132").Value = -11821.76

$ws.Range("H136").Value = 1355
$ws.Range("I136").Value = 912
$ws.Range("J136").Value = 1930.9
$ws.Range("K136").Value = 2736
$ws.Range("L136").Value = 5792.700000000001
$ws.Range("M136").Value = -186
$ws.Range("N136").Value = -10892.7

Write-Host "Applied all changes"